$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "A/C"
$ws.Range("C1").Value = 23
$ws.Range("D1").Value = $false

# Row 2
$ws.Range("A2").Value = "b"
$ws.Range("B2").Value = "Lâmpada"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = $false

# Row 3
$ws.Range("A3").Value = "c"
$ws.Range("B3").Value = "Televisor"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = $false
